$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet's displayed sheet name (workbook.xml sheet name)
$ws.Name = "2024 1 1"

# Update header cell B1: "Drink Name" -> "Create Date"
$ws.Range("B1").Value = "Create Date"

# Update data cells in row 2
$ws.Range("B2").Value = "2024-01-01T00:00"
$ws.Range("C2").Value = 2.0
$ws.Range("D2").Value = 65000.0

$wb.Save()
